# Add data for 2022-02-21: roll the "through" date in the sheet name and
# the February header forward from 02-12 to 02-13, and update the
# February / Total rows with the new cumulative counts.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Rename the sheet (workbook.xml <sheet name="..."/>)
$ws.Name = "Through 2022-02-13"

# Update the February row label (shared string)
$ws.Range("A3").Value = "February (through 02-13)"

# February row (row 3) updated counts - columns C:I (2016-2022); B (2015) & F (2019) unchanged
$ws.Range("C3").Value = 16
$ws.Range("D3").Value = 33
$ws.Range("E3").Value = 25
$ws.Range("G3").Value = 33
$ws.Range("H3").Value = 65
$ws.Range("I3").Value = 59

# Total row (row 4) updated counts - columns C:I (2016-2022); B (2015) & F (2019) unchanged
$ws.Range("C4").Value = 67
$ws.Range("D4").Value = 108
$ws.Range("E4").Value = 111
$ws.Range("G4").Value = 107
$ws.Range("H4").Value = 282
$ws.Range("I4").Value = 220
